# fix import authentication error!
# Rename sheet "totem" -> "totemConfig" and append two new rows of totem
# config data (rows 9 and 10) to the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet.
$ws.Name = "totemConfig"

# 2) Append the new data rows (columns A:L).
$newRows = @(
    @(7, 7, 1, $false, $false, 0, 73, 30, 3, 5, 1, 0),
    @(8, 8, 1, $false, $false, 0, 75, 30, 3, 5, 1, 0)
)

$startRow = 9
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $newRows[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}

# 3) Update the active selection to L9, matching the recorded edit state.
$ws.Range("L9").Select()
